$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now uses data originally from row 7
$ws.Range("D2").Value = 44181
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("O2").Value = "Provincia de Cautín"
$ws.Range("P2").Value = 560

# Row 3 now uses data originally from row 11
$ws.Range("D3").Value = 44186
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 14000
$ws.Range("O3").Value = "Región de La Araucanía"
$ws.Range("P3").Value = 560

# Row 4 now uses data originally from row 9
$ws.Range("D4").Value = 44159
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 8000
$ws.Range("O4").Value = "Región de La Araucanía"
$ws.Range("P4").Value = 320

# Row 5 now uses data originally from row 10
$ws.Range("D5").Value = 44159
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 8000
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 320

# Row 6 now uses data originally from row 21
$ws.Range("D6").Value = 44435
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 600

# Row 7 now uses data originally from row 8
$ws.Range("D7").Value = 44355
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 16000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 16000
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 640

# Row 8 now uses data originally from row 16
$ws.Range("D8").Value = 44434
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 600

# Row 9 now uses data originally from row 15
$ws.Range("D9").Value = 44452
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 13000
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 520

# Row 10 now uses data originally from row 20
$ws.Range("D10").Value = 44161
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("O10").Value = "Región de La Araucanía"
$ws.Range("P10").Value = 280

# Row 11 now uses data originally from row 13
$ws.Range("D11").Value = 44392
$ws.Range("J11").Value = 55
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17455
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 698

# Row 12 now uses data originally from row 23
$ws.Range("D12").Value = 44448
$ws.Range("J12").Value = 35
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 600

# Row 13 now uses data originally from row 5
$ws.Range("D13").Value = 44354
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 16000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 16000
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 640

# Row 14 now uses data originally from row 24
$ws.Range("D14").Value = 44167
$ws.Range("J14").Value = 95
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("O14").Value = "Región de La Araucanía"
$ws.Range("P14").Value = 280

# Row 15 now uses data originally from row 17
$ws.Range("D15").Value = 44210
$ws.Range("J15").Value = 110
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 16000
$ws.Range("M15").Value = 16000
$ws.Range("O15").Value = "Región de La Araucanía"
$ws.Range("P15").Value = 640

# Row 16 now uses data originally from row 18
$ws.Range("D16").Value = 44427
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 15000
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 600

# Row 17 now uses data originally from row 14
$ws.Range("D17").Value = 44168
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 7458
$ws.Range("O17").Value = "Región de La Araucanía"
$ws.Range("P17").Value = 298

# Row 18 now uses data originally from row 12
$ws.Range("D18").Value = 44371
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 15000
$ws.Range("O18").Value = "Provincia de Limarí"
$ws.Range("P18").Value = 600

# Row 19 now uses data originally from row 3
$ws.Range("D19").Value = 44356
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 14000
$ws.Range("O19").Value = "Provincia de Limarí"
$ws.Range("P19").Value = 560

# Row 20 now uses data originally from row 4
$ws.Range("D20").Value = 44160
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 8000
$ws.Range("O20").Value = "Región de La Araucanía"
$ws.Range("P20").Value = 320

# Row 21 now uses data originally from row 22
$ws.Range("D21").Value = 44175
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 8000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 8000
$ws.Range("O21").Value = "Región de La Araucanía"
$ws.Range("P21").Value = 320

# Row 22 now uses data originally from row 2
$ws.Range("D22").Value = 44357
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 14000
$ws.Range("L22").Value = 14000
$ws.Range("M22").Value = 14000
$ws.Range("O22").Value = "Provincia de Limarí"
$ws.Range("P22").Value = 560

# Row 23 now uses data originally from row 6
$ws.Range("D23").Value = 44162
$ws.Range("J23").Value = 260
$ws.Range("K23").Value = 7000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 7462
$ws.Range("O23").Value = "Región de La Araucanía"
$ws.Range("P23").Value = 298

# Row 24 now uses data originally from row 19
$ws.Range("D24").Value = 44176
$ws.Range("J24").Value = 20
$ws.Range("K24").Value = 11000
$ws.Range("L24").Value = 11000
$ws.Range("M24").Value = 11000
$ws.Range("O24").Value = "Región de La Araucanía"
$ws.Range("P24").Value = 440
